# Update the crypto price/volume snapshot (scraped GitHub Actions data refresh).
# Row 26/27 also swap Cosmos <-> Monero (source data changed ranking order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '39.428.46'
$ws.Range("E2").Value = '  +1.61%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.159.37'
$ws.Range("E3").Value = '  +3.15%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.00%  '

# Row 5: BNB
$ws.Range("D5").Value = '''227.59'
$ws.Range("E5").Value = '  -0.57%  '

# Row 6: XRP
$ws.Range("D6").Value = '''0.623'
$ws.Range("E6").Value = '  +0.89%  '

# Row 7: Solana
$ws.Range("D7").Value = '''64.27'
$ws.Range("E7").Value = '  +4.30%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  +2.61%  '

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.0860'
$ws.Range("E10").Value = '  +1.84%  '

# Row 11: TRON
$ws.Range("E11").Value = '  +0.43%  '

# Row 12: Chainlink
$ws.Range("D12").Value = '''16.03'
$ws.Range("E12").Value = '  +4.48%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '2.479.50'
$ws.Range("E13").Value = '  +3.16%  '

# Row 14: Avalanche
$ws.Range("D14").Value = '''22.28'
$ws.Range("E14").Value = '  +0.65%  '

# Row 15: Polygon
$ws.Range("D15").Value = '''0.812'
$ws.Range("E15").Value = '  +0.41%  '

# Row 16: Polkadot
$ws.Range("D16").Value = '''5.55'
$ws.Range("E16").Value = '  +1.15%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '2.158.04'
$ws.Range("E17").Value = '  +2.63%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '39.396.81'
$ws.Range("E18").Value = '  +1.62%  '

# Row 19: Litecoin
$ws.Range("D19").Value = '''71.78'
$ws.Range("E19").Value = '  -0.18%  '

# Row 20: Uniswap
$ws.Range("D20").Value = '''6.12'
$ws.Range("E20").Value = '  +0.81%  '

# Row 21: ShibaInu
$ws.Range("D21").Value = '0.0₃0855'
$ws.Range("E21").Value = '  +1.78%  '

# Row 22: BitcoinCash
$ws.Range("D22").Value = '''231.57'
$ws.Range("E22").Value = '  +1.46%  '

# Row 23: Dai
$ws.Range("E23").Value = '  +0.01%  '

# Row 24: Toncoin
$ws.Range("D24").Value = '''2.51'
$ws.Range("E24").Value = '  +6.38%  '

# Row 25: PancakeSwap
$ws.Range("D25").Value = '''2.36'
$ws.Range("E25").Value = '  +0.81%  '

# Row 26: Cosmos
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''172.26'
$ws.Range("E26").Value = '  +0.46%  '

# Row 27: Monero
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''9.52'
$ws.Range("E27").Value = '  -0.42%  '

# Row 28: Kaspa
$ws.Range("D28").Value = '''0.140'
$ws.Range("E28").Value = '  +1.78%  '

# Row 29: EthereumClassic
$ws.Range("D29").Value = '''19.87'
$ws.Range("E29").Value = '  +2.24%  '

# Row 30: ImmutableX
$ws.Range("D30").Value = '''1.41'
$ws.Range("E30").Value = '  -0.81%  '

# Row 31: WEMIXToken
$ws.Range("D31").Value = '''2.68'
$ws.Range("E31").Value = '  +6.52%  '

# Row 32: Stellar
$ws.Range("E32").Value = '  +0.59%  '

# Row 33: Filecoin
$ws.Range("D33").Value = '''4.62'
$ws.Range("E33").Value = '  +2.17%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").Value = '''4.75'
$ws.Range("E34").Value = '  -0.08%  '

# Row 35: THORChain
$ws.Range("D35").Value = '''7.09'
$ws.Range("E35").Value = '  +8.92%  '

# Row 36: Hedera
$ws.Range("D36").Value = '''0.0618'
$ws.Range("E36").Value = '  +0.08%  '

# Row 37: LidoDAOToken
$ws.Range("E37").Value = '  +0.43%  '

# Row 38: RenderToken
$ws.Range("E38").Value = '  -0.30%  '

# Row 39: BinanceUSD
$ws.Range("E39").Value = '  +0.16%  '

# Row 40: Aave
$ws.Range("D40").Value = '''104.16'
$ws.Range("E40").Value = '  +2.93%  '

# Row 41: VeChain
$ws.Range("E41").Value = '  +0.80%  '

# Row 42: InjectiveProtocol
$ws.Range("E42").Value = '  -2.26%  '

# Row 43: Maker
$ws.Range("D43").Value = '1.539.83'
$ws.Range("E43").Value = '  +0.39%  '

# Row 44: TrustWalletToken
$ws.Range("E44").Value = '  +3.60%  '

# Row 45: FraxShare
$ws.Range("D45").Value = '''7.89'
$ws.Range("E45").Value = '  +2.81%  '

# Row 46: HuobiToken
$ws.Range("E46").Value = '  +0.60%  '

# Row 47: Cronos
$ws.Range("D47").Value = '''0.0925'
$ws.Range("E47").Value = '  +1.53%  '

# Row 48: ARBITRUM
$ws.Range("E48").Value = '  +5.72%  '

# Row 49: FTXToken
$ws.Range("E49").Value = '  +2.86%  '

# Row 50: RocketPoolETH
$ws.Range("D50").Value = '2.361.87'
$ws.Range("E50").Value = '  +3.06%  '

# Row 51: MXToken
$ws.Range("D51").Value = '''2.97'
$ws.Range("E51").Value = '  +0.16%  '
